$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")

# Extend the monthly "PERIOD" date series in column A (rows 57-91) and
# fill in the newly-earned leave credit rows (60-65) with their EARNED
# value of 1.25, matching the continued monthly accrual pattern.

$ws.Range("A57").Value = 44957
$ws.Range("A58").Value = 44985
$ws.Range("A59").Value = 45016

$ws.Range("A60").Value = 45046
$ws.Range("C60").Value = 1.25

$ws.Range("A61").Value = 45077
$ws.Range("C61").Value = 1.25

$ws.Range("A62").Value = 45107
$ws.Range("C62").Value = 1.25

$ws.Range("A63").Value = 45138
$ws.Range("C63").Value = 1.25

$ws.Range("A64").Value = 45169
$ws.Range("C64").Value = 1.25

$ws.Range("A65").Value = 45199
$ws.Range("C65").Value = 1.25

$ws.Range("A66").Value = 45230
$ws.Range("A67").Value = 45260
$ws.Range("A68").Value = 45291
$ws.Range("A69").Value = 45322
$ws.Range("A70").Value = 45351
$ws.Range("A71").Value = 45382
$ws.Range("A72").Value = 45412
$ws.Range("A73").Value = 45443
$ws.Range("A74").Value = 45473
$ws.Range("A75").Value = 45504
$ws.Range("A76").Value = 45535
$ws.Range("A77").Value = 45565
$ws.Range("A78").Value = 45596
$ws.Range("A79").Value = 45626
$ws.Range("A80").Value = 45657
$ws.Range("A81").Value = 45688
$ws.Range("A82").Value = 45716
$ws.Range("A83").Value = 45747
$ws.Range("A84").Value = 45777
$ws.Range("A85").Value = 45808
$ws.Range("A86").Value = 45838
$ws.Range("A87").Value = 45869
$ws.Range("A88").Value = 45900
$ws.Range("A89").Value = 45930
$ws.Range("A90").Value = 45961
$ws.Range("A91").Value = 45991

$excel.Calculate()
